# Workflow.xlsx update — "All 2017 data K to Sap Plots"
# Adds the K-to-Sap dataset rows to the Expansion sheet and updates the
# active-sheet / selection state on both sheets.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")
$exp  = $wb.Worksheets.Item("Expansion")

# ---------------------------------------------------------------------
# Expansion sheet: new content (rows 3-7 gain cells, new rows 5-7 & 13)
# ---------------------------------------------------------------------

# Row 3 - tag the existing "Basic wrangling" row with the "All 2017 data" note
$exp.Range("E3").Value = "All 2017 data"

# Row 4 - "Convert k to sap" row gains a source link, a note, and the 2017 tag
$exp.Range("C4").Value = "https://github.com/AzuraLiu/MP.June/tree/main/Data/Processed/KtoSap"
$exp.Hyperlinks.Add($exp.Range("C4"), "https://github.com/AzuraLiu/MP.June/tree/main/Data/Processed/KtoSap")
$exp.Range("C4").Style = "Hyperlink"
$exp.Range("C4").WrapText = $true

$exp.Range("D4").Value = "Missing S3 March, S4 May, all Oct"
$exp.Range("E4").Value = "All 2017 data"

# Row 5 (new) - Sums
$exp.Range("A5").Value = "Sums"
$exp.Range("B5").Value = "Daily, daytime, & night sums"

$exp.Range("C5").Value = "https://github.com/AzuraLiu/MP.June/tree/main/Data/Processed/Sums/Raw"
$exp.Hyperlinks.Add($exp.Range("C5"), "https://github.com/AzuraLiu/MP.June/tree/main/Data/Processed/Sums/Raw")
$exp.Range("C5").Style = "Hyperlink"
$exp.Range("C5").WrapText = $true

$exp.Range("D5").Value = "Not cleaned"
$exp.Range("E5").Value = "All 2017 data"
$exp.Rows.Item(5).RowHeight = 28.8

# Row 6 (new) - Graph N vs S, by system (copied/adapted from Main)
$exp.Range("A6").Value = "Graph N vs S, by system"
$exp.Range("A6").WrapText = $true
$exp.Range("A6").VerticalAlignment = -4108
$exp.Range("A6").HorizontalAlignment = -4131

$exp.Range("B6").Value = "sort pairs, break by 10 days, graphed with VPD"
$exp.Range("B6").WrapText = $true
$exp.Range("B6").VerticalAlignment = -4108
$exp.Range("B6").HorizontalAlignment = -4131

$exp.Range("C6").Value = "https://github.com/AzuraLiu/MP.June/tree/main/Data/Plots/NS"
$exp.Hyperlinks.Add($exp.Range("C6"), "https://github.com/AzuraLiu/MP.June/tree/main/Data/Plots/NS")
$exp.Range("C6").Style = "Hyperlink"
$exp.Range("C6").WrapText = $true
$exp.Range("C6").VerticalAlignment = -4108

$exp.Range("D6").Value = "LBON17 missing Jan & Feb; HCON4 missing Dec"
$exp.Range("D6").WrapText = $true
$exp.Range("D6").VerticalAlignment = -4108

$exp.Range("E6").Value = "All 2017 Data"
$exp.Range("E6").WrapText = $true
$exp.Range("E6").VerticalAlignment = -4108

$exp.Rows.Item(6).RowHeight = 30

# Row 7 (new) - Cleaning Erorrs (copied/adapted from Main)
$exp.Range("A7").Value = "Cleaning Erorrs"
$exp.Range("A7").WrapText = $true
$exp.Range("A7").VerticalAlignment = -4108

$exp.Range("B7").Value = "Compare VPD, N<S?, break starts, create new file with cleaned dat"
$exp.Range("B7").WrapText = $true
$exp.Range("B7").VerticalAlignment = -4108

$exp.Range("D7").Value = "In progress"
$exp.Range("D7").WrapText = $true
$exp.Range("D7").VerticalAlignment = -4108

# Row 13 (new) - trailing note
$exp.Range("D13").Value = " "

# ---------------------------------------------------------------------
# Selection / active sheet updates
# ---------------------------------------------------------------------

# Expansion is no longer the tab shown on re-open; Main becomes active.
$exp.Range("A2").Select()

$main.Activate()
$main.Range("K13").Select()

Write-Output "Applied K-to-Sap Plots update"
